$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 617, shifting existing rows 617-679 down to 618-680
$ws.Rows.Item(617).Insert()

# Fill the new row 617 with the new weekly entry values
$ws.Cells.Item(617, 1).Value = 3
$ws.Cells.Item(617, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(617, 3).Value = "Coquimbo"
$ws.Cells.Item(617, 4).Value = 45265
$ws.Cells.Item(617, 5).Value = 5
$ws.Cells.Item(617, 6).Value = 100112027
$ws.Cells.Item(617, 7).Value = "Melón"
$ws.Cells.Item(617, 8).Value = "Tuna"
$ws.Cells.Item(617, 9).Value = "Primera"
$ws.Cells.Item(617, 10).Value = 130
$ws.Cells.Item(617, 11).Value = 2000
$ws.Cells.Item(617, 12).Value = 2000
$ws.Cells.Item(617, 13).Value = 2000
$ws.Cells.Item(617, 14).Value = "$/unidad"
$ws.Cells.Item(617, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(617, 16).Value = 2000
$ws.Cells.Item(617, 17).Value = 1
$ws.Cells.Item(617, 18).Value = "Hortaliza"
